$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "42.249.66"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "2.177.04"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "252.54"
$ws.Range("E5").Value = "  +5.95%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "73.43"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "40.78"
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0910"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.76"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.101"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.504.43"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "14.21"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.179.15"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "0.762"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "42.149.73"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0000102"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "70.64"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.87"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "226.57"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "9.49"
$ws.Range("E24").Value = "  -5.37%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "10.46"
$ws.Range("E26").Value = "  -2.76%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "3.32"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  +6.68%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "36.79"
$ws.Range("E30").Value = "  +11.73%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "169.18"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "20.01"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0805"
$ws.Range("E33").Value = "  +3.47%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.11"
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.120"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.107"
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.23"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.0332"
$ws.Range("E38").Value = "  +6.11%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "11.91"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.05"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.195"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "59.06"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "102.27"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "2.46"
$ws.Range("E45").Value = "  +12.13%  "
$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D46").Value = "0.469"
$ws.Range("E46").Value = "  +15.01%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "8.29"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0968"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "1.09"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.12"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  +0.77%  "
